$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.022.82'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '2.415.15'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''562.71'
$ws.Range('E5').Value = '  +1.27%  '
$ws.Range('D6').Value = '''142.52'
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''0.531'
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('E9').Value = '  +0.44%  '
$ws.Range('E10').Value = '  -1.80%  '
$ws.Range('D11').Value = '''5.22'
$ws.Range('E11').Value = '  -3.30%  '
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('D13').Value = '''25.58'
$ws.Range('E13').Value = '  -2.23%  '
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('D15').Value = '2.854.92'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').Value = '61.901.93'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').Value = '2.412.17'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').Value = '''11.23'
$ws.Range('D19').Value = '''321.41'
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('D20').Value = '''6.82'
$ws.Range('E20').Value = '  +1.90%  '
$ws.Range('D21').Value = '''4.13'
$ws.Range('E21').Value = '  -1.38%  '
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('E23').Value = '  +1.94%  '
$ws.Range('D24').Value = '''1.75'
$ws.Range('E24').Value = '  -0.57%  '
$ws.Range('D25').Value = '''8.63'
$ws.Range('E25').Value = '  -5.22%  '
$ws.Range('D26').Value = '''564.19'
$ws.Range('E26').Value = '  -2.35%  '
$ws.Range('D27').Value = '2.533.90'
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('D29').Value = '0.0₃0932'
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('E30').Value = '  -2.06%  '
$ws.Range('E31').Value = '  -3.92%  '
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('E33').Value = '  +1.24%  '
$ws.Range('E34').Value = '  -2.41%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').Value = '''4.78'
$ws.Range('E36').Value = '  -0.85%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '''153.14'
$ws.Range('E37').Value = '  +3.29%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').Value = '''5.43'
$ws.Range('E38').Value = '  -3.94%  '
$ws.Range('D39').Value = '''0.379'
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('D40').Value = '''18.55'
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('D41').Value = '''1.80'
$ws.Range('E41').Value = '  -3.66%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = '''149.00'
$ws.Range('E43').Value = '  -1.29%  '
$ws.Range('E44').Value = '  -2.28%  '
$ws.Range('D45').Value = '''3.61'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('D46').Value = '''0.0531'
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '''0.594'
$ws.Range('E47').Value = '  +1.12%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '''19.82'
$ws.Range('E48').Value = '  -2.47%  '
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('E51').Value = '  +0.67%  '
